$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3.53
$ws.Range("D2").Value = 3.58
$ws.Range("E2").Value = 3.66
$ws.Range("F2").Value = 3.69

# Row 3
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 4.04
$ws.Range("D3").Value = 4.04
$ws.Range("F3").Value = 4.13

# Row 4
$ws.Range("B4").Value = 3.83
$ws.Range("D4").Value = 3.83
$ws.Range("E4").Value = 3.92
$ws.Range("F4").Value = 3.76

# Row 5
$ws.Range("B5").Value = 3.51
$ws.Range("C5").Value = 3.47
$ws.Range("D5").Value = 3.48
$ws.Range("E5").Value = 3.39
$ws.Range("F5").Value = 3.16

# Row 6
$ws.Range("B6").Value = 3.18
$ws.Range("C6").Value = 3.29
$ws.Range("D6").Value = 3.2
$ws.Range("E6").Value = 3.17
$ws.Range("F6").Value = 3.35

# Row 7
$ws.Range("B7").Value = 2.94
$ws.Range("C7").Value = 3.04
$ws.Range("E7").Value = 2.88
$ws.Range("F7").Value = 2.81

# Row 8
$ws.Range("B8").Value = 4.53
$ws.Range("E8").Value = 4.52
$ws.Range("F8").Value = 4.53
